# Bug fix: connector/line shapes on the slide were created with an
# essentially invisible line weight (1 EMU, i.e. ~0.00008 pt). Restore
# the intended hairline width of 1pt (12700 EMU) for every straight-line
# connector shape in the deck.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        # msoLine (straight connector shapes using prstGeom "line")
        if ($shp.Type -eq 9) {
            if ([Math]::Abs($shp.Line.Weight - (1 / 12700)) -lt 0.00001) {
                $shp.Line.Weight = 1
            }
        }
    }
}
